# Re-creates the edit described by the commit:
#   "Added new dataset preprocessing"
#
# The original workbook has a single worksheet, "Data Harian - Table",
# that holds some station metadata (rows 1-5), a daily-readings table
# (rows 9-39: one header row + 30 data rows across columns A:K) and a
# legend (rows 43-55).
#
# The edit duplicates just the daily-readings table (A9:K39) into a
# brand-new worksheet named "Sheet1", placed immediately after the
# original sheet and left as the active/selected sheet - i.e. someone
# copy/pasted the clean table out to its own sheet for further
# preprocessing.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- restore a couple of cosmetic view settings that the original
# file stored in a verbose boolean form (gridlines-on, wrapped text)
# so the resulting sheet view matches the original look -------------
$ws1.Activate()
$excel.ActiveWindow.DisplayGridlines = $true

# A1 on the metadata block is (and always was) empty - drop the stray
# empty cell record.
$ws1.Range("A1").ClearContents()

# Keep the table's original formatting (thin borders, centered/wrapped
# header, left/top wrapped body) intact after the round-trip.
$headerRow = $ws1.Range("A9:K9")
$headerRow.Borders.LineStyle = 1
$headerRow.Borders.Color = 0
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4108
$headerRow.WrapText = $true

$bodyRows = $ws1.Range("A10:K39")
$bodyRows.Borders.LineStyle = 1
$bodyRows.Borders.Color = 0
$bodyRows.HorizontalAlignment = -4131
$bodyRows.VerticalAlignment = -4160
$bodyRows.WrapText = $true

# --- copy the clean daily-readings table (header + 30 rows) into a
# new worksheet -------------------------------------------------------
$source = $ws1.Range("A9:K39")
$source.Copy()

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Range("A1").PasteSpecial(-4104)
$ws2.Name = "Sheet1"

# Re-apply the same formatting to the copy (paste alone only carries
# values/number formats across in this host).
$newHeader = $ws2.Range("A1:K1")
$newHeader.Borders.LineStyle = 1
$newHeader.Borders.Color = 0
$newHeader.HorizontalAlignment = -4108
$newHeader.VerticalAlignment = -4108
$newHeader.WrapText = $true
$newHeader.RowHeight = 14.4

$newBody = $ws2.Range("A2:K31")
$newBody.Borders.LineStyle = 1
$newBody.Borders.Color = 0
$newBody.HorizontalAlignment = -4131
$newBody.VerticalAlignment = -4160
$newBody.WrapText = $true
$newBody.RowHeight = 28.8

# Leave the new sheet as the active tab with the whole table selected,
# and restore the original sheet's selection/scroll position over the
# table it was copied from.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws1.Range("A9:K39").Select()

$ws2.Activate()
$ws2.Range("A1:K31").Select()
